# Wellness.xlsx - "Mise a jour de l'application"
# Append 10 new daily wellness entries (rows 410-419) for the 2025-09-26
# session (Excel date serial 45926), following the same column layout as
# the existing data: Date, Nom du joueur, Volume, Intensite, Fatigue,
# Douleur, Localisation douleur, Plaisir, Charge(=Volume*Intensite).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$nbsp = [char]0x00A0

# Name, Volume(C), Intensite(D), Fatigue(E), Douleur(F), Localisation
# douleur(G) and Plaisir(H) for each new row, in sheet order.
$rows = @(
    @{ Name = "Jeremie Laurent";  C = 60; D = 6; E = 6; F = 0; G = $null;             H = 6  },
    @{ Name = "Emmanuel Valey";   C = 60; D = 2; E = 5; F = 2; G = "Ischio" + $nbsp;  H = 7  },
    @{ Name = "Ilan Ihaddadene";  C = 60; D = 6; E = 6; F = 0; G = $null;             H = 10 },
    @{ Name = "Mattheo Haon";     C = 60; D = 5; E = 7; F = 2; G = "Cheville";        H = 7  },
    @{ Name = "Omar Benyounes";   C = 60; D = 4; E = 5; F = 0; G = $null;             H = 6  },
    @{ Name = "Hedi Nasri";       C = 60; D = 3; E = 3; F = 4; G = "Dos";             H = 7  },
    @{ Name = "Naim Ighbane";     C = 60; D = 5; E = 4; F = 3; G = "Cheville";        H = 4  },
    @{ Name = "Karim Belmahi";    C = 60; D = 6; E = 3; F = 0; G = $null;             H = 10 },
    @{ Name = "Sofiane Belle";    C = 60; D = 2; E = 3; F = 0; G = $null;             H = 0  },
    @{ Name = "Amir Etien";       C = 60; D = 6; E = 6; F = 0; G = $null;             H = 5  }
)

$startRow = 410
$endRow = $startRow + $rows.Count - 1
$dateSerial = 45926

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $d = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $d.Name
    $ws.Cells.Item($r, 3).Value = $d.C
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E
    $ws.Cells.Item($r, 6).Value = $d.F
    if ($d.G) {
        $ws.Cells.Item($r, 7).Value = $d.G
    }
    $ws.Cells.Item($r, 8).Value = $d.H
}

# A: date, formatted like the rows above it (built-in short-date style)
$ws.Range("A$startRow`:A$endRow").NumberFormat = "m/d/yy"

# B, C-F, H: regular black "Helvetica Neue" 10pt style used throughout
# the sheet's data rows.
$dataFont = $ws.Range("B$startRow`:F$endRow").Font
$dataFont.Color = 0
$dataFont.Name = "Helvetica Neue"
$dataFont.Size = 10

$hFont = $ws.Range("H$startRow`:H$endRow").Font
$hFont.Color = 0
$hFont.Name = "Helvetica Neue"
$hFont.Size = 10

# G: Localisation douleur - empty and filled cells alike use the
# "Helvetica" 12pt (theme colour) style used for the rest of the column.
$gFont = $ws.Range("G$startRow`:G$endRow").Font
$gFont.Name = "Helvetica"
$gFont.Size = 12

# I: Charge = Volume * Intensite, written as one pass so it becomes a
# single shared formula across the new block (matching the existing
# C#*D# pattern used for every row above).
$ws.Range("I$startRow`:I$endRow").Formula = "=C$startRow*D$startRow"

# Reflect the same scroll position / selection the author ended up with
# after adding the new rows.
$ws.Application.ActiveWindow.ScrollRow = 390
$ws.Range("K416").Select()
